# Add "PayGradesTest" rows (Class 4 / Grade A / Grade B / Grade C) and a
# "Duplicate" column to the existing Pay Grades data block on the
# "Test Data" sheet (sheet 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Extend the existing header/first-data row formatting out to column G
#     (a "Duplicate" column is being inserted before the trailing "Runmode"
#     column, just like the other test-data blocks on this sheet already
#     have) ---
$ws.Range("A34:F34").Copy() | Out-Null
$ws.Range("A34:G34").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A35:F35").Copy() | Out-Null
$ws.Range("A35:G35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header row: insert "Duplicate" before the trailing "Runmode" column.
$ws.Range("F34").Value = "Duplicate"
$ws.Range("G34").Value = "Runmode"

# Existing "Class 3" row: minimum salary becomes the numeral "10000"
# (replacing the spelled-out "ten thousand") and the new Duplicate column
# is "N".
$ws.Range("D35").Value = "'10000"
$ws.Range("G35").Value = "N"

# --- Add five new Pay Grade data rows (36-40), reusing the formatting of
#     row 35 ---
$ws.Range("A35:G35").Copy() | Out-Null
$ws.Range("A36:G40").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A36").Value = "Chrome"
$ws.Range("B36").Value = "Class 4"
$ws.Range("C36").Value = "USD"
$ws.Range("D36").Value = "'10000"
$ws.Range("E36").Value = "'100000"
$ws.Range("F36").Value = "N"
$ws.Range("G36").Value = "N"

$ws.Range("A37").Value = "Chrome"
$ws.Range("B37").Value = "Grade A"
$ws.Range("C37").Value = "AUD"
$ws.Range("D37").Value = "'350000"
$ws.Range("E37").Value = "'250000"
$ws.Range("F37").Value = "N"
$ws.Range("G37").Value = "Y"

$ws.Range("A38").Value = "Chrome"
$ws.Range("B38").Value = "Grade B"
$ws.Range("C38").Value = "AUD"
$ws.Range("D38").Value = "'350000"
$ws.Range("E38").Value = "'2500000"
$ws.Range("F38").Value = "N"
$ws.Range("G38").Value = "N"

$ws.Range("A39").Value = "Chrome"
$ws.Range("C39").Value = "SGD"
$ws.Range("D39").Value = "'350000"
$ws.Range("E39").Value = "'2500000"
$ws.Range("F39").Value = "N"
$ws.Range("G39").Value = "N"

$ws.Range("A40").Value = "Chrome"
$ws.Range("B40").Value = "Grade C"
$ws.Range("F40").Value = "N"
$ws.Range("G40").Value = "Y"

Write-Output "Pay grade rows added"
